$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ITEMS")
$ws.Range("A1").Value = "TEST"
